$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "trade_date" column (I) for all existing data rows (2-356) switches
# from a date-only number format to a date+time number format (matching
# the "datetime" column B), which re-points those cells at cellXfs style 2
# instead of style 3.
$ws.Range("I2:I356").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Four new trading-day rows appended at the bottom of the sheet.
# Columns: close, datetime, exchange_code, high, low, open, stock_code,
#          volume, trade_date, isin, company, symbol, source
$newRows = @(
    @(166.07, 46049, "NSE", 169.04, 157,    160,    "DCXSYS", 1067707, 46049, "INE0KL801015", "DCX System Ltd", "DCXSYS", "BREEZE"),
    @(177.23, 46050, "NSE", 178.16, 168.52, 168.52, "DCXSYS", 1089999, 46050, "INE0KL801015", "DCX System Ltd", "DCXSYS", "BREEZE"),
    @(171.92, 46051, "NSE", 179.23, 171.5,  177.61, "DCXSYS", 535071,  46051, "INE0KL801015", "DCX System Ltd", "DCXSYS", "BREEZE"),
    @(177.09, 46052, "NSE", 181.95, 167.92, 170.93, "DCXSYS", 904390,  46052, "INE0KL801015", "DCX System Ltd", "DCXSYS", "BREEZE")
)

$r = 357
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]

    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 9).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]

    $r = $r + 1
}
